$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for the Price column cells whose new values would
# otherwise be auto-converted to numbers (losing formatting like trailing zeros)
# before assigning their values, matching the original inline-string text.

$ws.Range('D2').Value = '88.532.80'
$ws.Range('E2').Value = '  +8.72%  '
$ws.Range('D3').Value = '3.335.72'
$ws.Range('E3').Value = '  +5.03%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '220.03'
$ws.Range('E5').Value = '  +6.01%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '650.44'
$ws.Range('E6').Value = '  +2.71%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.378'
$ws.Range('E7').Value = '  +28.97%  '
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('E9').Value = '  +2.71%  '
$ws.Range('D10').Value = '3.334.78'
$ws.Range('E10').Value = '  +4.98%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.587'
$ws.Range('E11').Value = '  -0.54%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000282'
$ws.Range('E12').Value = '  +7.92%  '
$ws.Range('B13').Value = 'TRON'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.168'
$ws.Range('E13').Value = '  +1.75%  '
$ws.Range('B14').Value = 'Avalanche'
$ws.Range('C14').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '35.49'
$ws.Range('E14').Value = '  +11.27%  '
$ws.Range('D15').Value = '3.946.54'
$ws.Range('E15').Value = '  +4.77%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.50'
$ws.Range('E16').Value = '  +2.59%  '
$ws.Range('D17').Value = '88.117.52'
$ws.Range('E17').Value = '  +8.08%  '
$ws.Range('D18').Value = '3.317.05'
$ws.Range('E18').Value = '  +4.35%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '14.68'
$ws.Range('E19').Value = '  +3.87%  '
$ws.Range('B20').Value = 'SuiNetwork'
$ws.Range('C20').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '3.15'
$ws.Range('E20').Value = '  -2.93%  '
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.72'
$ws.Range('E21').Value = '  +5.66%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '457.16'
$ws.Range('E22').Value = '  +4.14%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.49'
$ws.Range('E23').Value = '  +6.09%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.61'
$ws.Range('E24').Value = '  +9.91%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '12.67'
$ws.Range('E25').Value = '  +12.48%  '
$ws.Range('D26').Value = '3.509.13'
$ws.Range('E26').Value = '  +4.73%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '78.71'
$ws.Range('E27').Value = '  +2.32%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.201'
$ws.Range('E28').Value = '  +44.84%  '
$ws.Range('E29').Value = '  +1.89%  '
$ws.Range('E30').Value = '  -0.28%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '607.60'
$ws.Range('E31').Value = '  +7.32%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '9.44'
$ws.Range('E32').Value = '  +3.37%  '
$ws.Range('E33').Value = '  +7.73%  '
$ws.Range('E34').Value = '  -0.23%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.11'
$ws.Range('E35').Value = '  +3.86%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '7.28'
$ws.Range('E36').Value = '  +22.71%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.146'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '23.26'
$ws.Range('E38').Value = '  +0.42%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.19'
$ws.Range('E39').Value = '  +8.18%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.422'
$ws.Range('E40').Value = '  +1.97%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '21.85'
$ws.Range('E41').Value = '  +5.21%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.00'
$ws.Range('E42').Value = '  -0.03%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.08'
$ws.Range('E43').Value = '  +0.40%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '160.07'
$ws.Range('E44').Value = '  +0.04%  '
$ws.Range('B45').Value = 'USDe'
$ws.Range('C45').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.00'
$ws.Range('E45').Value = '  +0.00%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '191.31'
$ws.Range('E46').Value = '  +1.25%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.44'
$ws.Range('E47').Value = '  +7.31%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '46.43'
$ws.Range('E48').Value = '  +4.79%  '
$ws.Range('B49').Value = 'Mantle'
$ws.Range('C49').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.789'
$ws.Range('E49').Value = '  +0.11%  '
$ws.Range('B50').Value = 'Filecoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '4.43'
$ws.Range('E50').Value = '  +3.72%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.663'
$ws.Range('E51').Value = '  +3.86%  '
